# Update countries & provincias Spain
# Applies the COVID data refresh described by the commit:
#  - Re-sorted leaderboard rows: Emiratos Arabes Unidos overtakes Dinamarca,
#    Uganda overtakes Haiti & Polinesia Francesa, Puerto Rico/Eritrea swap.
#  - Refreshed case counts for Suiza, Kazajistan, Barein and the above.
#  - Updated the "last refreshed" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 13:52"

# Helper-free, explicit per-row update: country name (col A) + B..H values
function Set-Row {
    param($row, $name, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 18: Suiza - refreshed counts, same rank
Set-Row 18 "Suiza" 28063 119 18600 8027 386 7 1436

# Rows 36/37: Emiratos Arabes Unidos overtakes Dinamarca
Set-Row 36 "Emiratos Arabes Unidos" 7755 490 1443 6266 1 3 46
Set-Row 37 "Dinamarca" 7695 180 4312 3019 84 0 364

# Row 64: Kazajistan - refreshed counts, same rank
Set-Row 64 "Kazajistan" 1967 115 476 1472 22 0 19

# Row 65: Barein - refreshed counts, same rank
Set-Row 65 "Barein" 1952 45 783 1162 2 0 7

# Rows 158/159/160: Uganda overtakes Haiti and Polinesia Francesa
Set-Row 158 "Uganda" 58 2 38 20 0 0 0
Set-Row 159 "Haiti" 57 0 0 54 0 0 3
Set-Row 160 "Polinesia Francesa" 56 0 19 37 1 0 0

# Rows 166/167: Puerto Rico and Eritrea swap order
Set-Row 166 "Puerto Rico" 39 0 1 36 0 0 2
Set-Row 167 "Eritrea" 39 0 3 36 0 0 0
